# "Finished Week 13 logging" — update Rushing & Receiving stat tables with
# this week's numbers, renumber the trailing week index for players who
# slot in after the newly logged player, and append G.Kittle's Rushing row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# J.Garoppolo: 3DATT 9 -> 10
$rushing.Range("E2").Value = 10

# E.Mitchell: updated week totals
$rushing.Range("C5").Value = 92
$rushing.Range("D5").Value = 52
$rushing.Range("E5").Value = 9
$rushing.Range("F5").Value = 16

# Week-index renumbering (a week slot was filled in, shifting these down by one)
$rushing.Range("A9").Value = 7
$rushing.Range("A10").Value = 8

# New row: G.Kittle's rushing stats for the week
$rushing.Range("A10").Copy($rushing.Range("A11"))
$rushing.Range("A11").Value = 9
$rushing.Range("B11").Value = "G.Kittle"
$rushing.Range("C11").Value = 0
$rushing.Range("D11").Value = 1
$rushing.Range("E11").Value = 0
$rushing.Range("F11").Value = 0

# ---------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# E.Mitchell: updated week totals
$receiving.Range("C3").Value = 16
$receiving.Range("D3").Value = 15

# K.Juszczyk: updated week totals
$receiving.Range("C5").Value = 20
$receiving.Range("D5").Value = 19
$receiving.Range("E5").Value = 3

# B.Aiyuk: week index + updated week totals
$receiving.Range("A7").Value = 5
$receiving.Range("C7").Value = 39
$receiving.Range("D7").Value = 26
$receiving.Range("E7").Value = 15
$receiving.Range("F7").Value = 9

# M.Sanu: week index only
$receiving.Range("A8").Value = 6

# T.Sherfield: week index + updated week totals
$receiving.Range("A9").Value = 7
$receiving.Range("C9").Value = 14
$receiving.Range("D9").Value = 7
$receiving.Range("G9").Value = 3

# J.Jennings: week index + updated week totals
$receiving.Range("A10").Value = 8
$receiving.Range("C10").Value = 11
$receiving.Range("D10").Value = 7

# G.Kittle: week index + updated week totals
$receiving.Range("A11").Value = 9
$receiving.Range("C11").Value = 65
$receiving.Range("D11").Value = 52
$receiving.Range("E11").Value = 21
$receiving.Range("F11").Value = 15

# R.Dwelley: week index only
$receiving.Range("A12").Value = 10

# C.Woerner: week index only
$receiving.Range("A13").Value = 11

# Leave the workbook with the Rushing tab active/selected, matching where
# logging finished up.
$rushing.Activate()
